$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet3 = "Translated_Sheet1" (was A1:B5, becomes A1:A5)
# Sheet4 = "Translated_Sheet2" (was A1:B4, becomes A1:A4)
# Column B ("translated_magyar oszlop" / python-list-repr values) is dropped
# in both sheets; column A keeps the Hungarian header + gets the plain
# (un-bracketed) English translation text. Formatting is copied from the
# already-correctly-formatted Sheet1 / Sheet2 (which show the same pattern
# for the Hungarian-language sheets) so the exact same style/border/fill
# records get reused.
# ---------------------------------------------------------------------------

$wsSrc1 = $wb.Worksheets.Item("Sheet1")
$wsSrc2 = $wb.Worksheets.Item("Sheet2")
$ws1 = $wb.Worksheets.Item("Translated_Sheet1")
$ws2 = $wb.Worksheets.Item("Translated_Sheet2")

# ---------------------------------------------------------------------------
# Translated_Sheet1
# ---------------------------------------------------------------------------

# Drop the second (translated_*) column entirely -> dimension becomes A1:A5
$ws1.Columns.Item(2).Delete()

# New plain-text English content
$ws1.Range("A1").Value = "Hungarian column"
$ws1.Range("A2").Value = "Apples"
$ws1.Range("A3").Value = "I don't think that's going to be difficult."
$ws1.Range("A4").Value = "We'll see"
$ws1.Range("A5").Value = "Other vehicles"

# A1 loses the bold/boxed/centered header formatting it used to share with B1
$ws1.Range("A1").ClearFormats()

# Reuse the exact cell formatting already present on Sheet1 (same layout,
# same highlight / border / font-color scheme) by copying formats only
$wsSrc1.Range("A2").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$wsSrc1.Range("A3").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

$wsSrc1.Range("A4").Copy()
$ws1.Range("A4").PasteSpecial(-4122)

$wsSrc1.Range("A5").Copy()
$ws1.Range("A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Matching row heights for the thick-bordered rows
$ws1.Rows.Item(4).RowHeight = 15
$ws1.Rows.Item(5).RowHeight = 15

# Column A width (best-fit sized to the longest translated string)
$ws1.Columns.Item(1).ColumnWidth = 31.1666666666667

# Page margins / setup
$ws1.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$ws1.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$ws1.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$ws1.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$ws1.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Translated_Sheet2
# ---------------------------------------------------------------------------

# Drop the second (translated_*) column entirely -> dimension becomes A1:A4
$ws2.Columns.Item(2).Delete()

# New plain-text English content
$ws2.Range("A1").Value = "Hungarian column"
$ws2.Range("A2").Value = "the roasting"
$ws2.Range("A3").Value = "Other, of a kind used for the manufacture of goods"
$ws2.Range("A4").Value = "That's a whole sentence."

# A1 loses the bold/boxed/centered header formatting it used to share with B1
$ws2.Range("A1").ClearFormats()

# Column A width (best-fit sized to the longest translated string)
$ws2.Columns.Item(1).ColumnWidth = 17.1666666666667

# Page margins (no pageSetup orientation override on this sheet)
$ws2.PageSetup.LeftMargin = $excel.InchesToPoints(0.7)
$ws2.PageSetup.RightMargin = $excel.InchesToPoints(0.7)
$ws2.PageSetup.TopMargin = $excel.InchesToPoints(0.75)
$ws2.PageSetup.BottomMargin = $excel.InchesToPoints(0.75)
$ws2.PageSetup.HeaderMargin = $excel.InchesToPoints(0.3)
$ws2.PageSetup.FooterMargin = $excel.InchesToPoints(0.3)

Write-Host "edits applied"
